# Generate Report for Handoff
#
# Swaps the display order of the two localized files
# (9928514f-...md and 11fda2fc-...md) on every sheet: the row that used to
# show "9928514f..." first now shows "11fda2fc..." first (row 2) and
# "9928514f..." second (row 3) - along with their related handoff-package
# hyperlinks - and refreshes the "Latest Handoff Datetime" for the new
# handoff run.

$wb = $excel.ActiveWorkbook

function Set-HL {
    param($ws, $cellRef, $displayText, $address)
    $ws.Range($cellRef).Value = $displayText
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $cellRef) {
            $hl.Address = $address
            $hl.TextToDisplay = $displayText
        }
    }
}

$md9928514f = "9928514f-2900-4a9c-9c54-d3e9e080dbf9.md"
$md11fda2fc = "11fda2fc-3867-49e3-bd1b-8de5d014f76a.md"

$mdUrl9928514f = "https://github.com/OpenLocalizationTest/oltest/blob/46f4beb0e240c8f1b09ff54ad46202877d6c4187/e2e/9928514f-2900-4a9c-9c54-d3e9e080dbf9.md"
$mdUrl11fda2fc = "https://github.com/OpenLocalizationTest/oltest/blob/46f4beb0e240c8f1b09ff54ad46202877d6c4187/e2e/11fda2fc-3867-49e3-bd1b-8de5d014f76a.md"

# ---- Sheet "Overview" ----
$ws = $wb.Worksheets.Item("Overview")
Set-HL $ws '$A$2' $md11fda2fc $mdUrl11fda2fc
Set-HL $ws '$A$3' $md9928514f $mdUrl9928514f

# ---- Sheet "zh-cn" ----
$ws = $wb.Worksheets.Item("zh-cn")

$xlf9928514fZhCn = "9928514f-2900-4a9c-9c54-d3e9e080dbf9.3e61cd54855f4dd589534a2d2eb0bc90365b634d.zh-cn.xlf"
$xlf11fda2fcZhCn = "11fda2fc-3867-49e3-bd1b-8de5d014f76a.ab49440f143130ce0ee6ea0532637013b9fe8bec.zh-cn.xlf"

$xlfUrl9928514fZhCn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dedec3aae0f29191a7f7d96e15e93cb6db4b43e7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/9928514f-2900-4a9c-9c54-d3e9e080dbf9.3e61cd54855f4dd589534a2d2eb0bc90365b634d.zh-cn.xlf"
$xlfUrl11fda2fcZhCn = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dedec3aae0f29191a7f7d96e15e93cb6db4b43e7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/11fda2fc-3867-49e3-bd1b-8de5d014f76a.ab49440f143130ce0ee6ea0532637013b9fe8bec.zh-cn.xlf"

Set-HL $ws '$A$2' $md11fda2fc $mdUrl11fda2fc
Set-HL $ws '$C$2' $xlf11fda2fcZhCn $xlfUrl11fda2fcZhCn
Set-HL $ws '$A$3' $md9928514f $mdUrl9928514f
Set-HL $ws '$C$3' $xlf9928514fZhCn $xlfUrl9928514fZhCn

# New handoff datetime for this run (shared by both rows, as before)
$ws.Range("D2").Value = "2016-03-03 08:33:45"
$ws.Range("D3").Value = "2016-03-03 08:33:45"

# ---- Sheet "de-de" ----
$ws = $wb.Worksheets.Item("de-de")

$xlf9928514fDeDe = "9928514f-2900-4a9c-9c54-d3e9e080dbf9.3e61cd54855f4dd589534a2d2eb0bc90365b634d.de-de.xlf"
$xlf11fda2fcDeDe = "11fda2fc-3867-49e3-bd1b-8de5d014f76a.ab49440f143130ce0ee6ea0532637013b9fe8bec.de-de.xlf"

$xlfUrl9928514fDeDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2b6c5b3f172d607ba340206a3442fc830c6b75e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/9928514f-2900-4a9c-9c54-d3e9e080dbf9.3e61cd54855f4dd589534a2d2eb0bc90365b634d.de-de.xlf"
$xlfUrl11fda2fcDeDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c2b6c5b3f172d607ba340206a3442fc830c6b75e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/11fda2fc-3867-49e3-bd1b-8de5d014f76a.ab49440f143130ce0ee6ea0532637013b9fe8bec.de-de.xlf"

Set-HL $ws '$A$2' $md11fda2fc $mdUrl11fda2fc
Set-HL $ws '$C$2' $xlf11fda2fcDeDe $xlfUrl11fda2fcDeDe
Set-HL $ws '$A$3' $md9928514f $mdUrl9928514f
Set-HL $ws '$C$3' $xlf9928514fDeDe $xlfUrl9928514fDeDe

# New handoff datetime for this run (shared by both rows, as before)
$ws.Range("D2").Value = "2016-03-03 08:33:57"
$ws.Range("D3").Value = "2016-03-03 08:33:57"
